$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts everything right)
$ws.Columns("A").Insert()

# Header for new column
$ws.Range("A1").Value = "Fonte"

# Fill rows 2-49 (old data rows) with "CORE01"
$ws.Range("A2:A49").Value = "CORE01"

# Update selection to mimic the authored edit
$ws.Range("A3:A49").Select()
